# Refresh the cryptos price/volume table (GitHub Actions scheduled update).
# Column D ("Price") is always stored as literal text in this sheet (even
# when the value looks numeric, e.g. "536.67"), so every Price cell we touch
# is written with a leading apostrophe. That stops Excel from re-parsing it
# as a number, which would silently normalise the text (drop trailing zeros,
# switch to scientific notation, etc.) and break the Price column's type.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'" + '61.246.90'
$ws.Range("E2").Value = '  -4.25%  '
$ws.Range("D3").Value = "'" + '2.991.71'
$ws.Range("E3").Value = '  -3.35%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").Value = "'" + '535.17'
$ws.Range("E5").Value = '  -1.29%  '
$ws.Range("D6").Value = "'" + '134.75'
$ws.Range("E6").Value = '  -2.07%  '
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("D8").Value = "'" + '2.985.42'
$ws.Range("E8").Value = '  -3.27%  '
$ws.Range("E9").Value = '  -0.51%  '
$ws.Range("B10").Value = 'Toncoin'
$ws.Range("C10").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D10").Value = "'" + '6.15'
$ws.Range("E10").Value = '  -0.33%  '
$ws.Range("B11").Value = 'Dogecoin'
$ws.Range("C11").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D11").Value = "'" + '0.148'
$ws.Range("E11").Value = '  -5.55%  '
$ws.Range("D12").Value = "'" + '0.447'
$ws.Range("E12").Value = '  -2.92%  '
$ws.Range("E13").Value = '  -2.52%  '
$ws.Range("D14").Value = "'" + '33.94'
$ws.Range("D15").Value = "'" + '3.477.51'
$ws.Range("E15").Value = '  -3.45%  '
$ws.Range("E16").Value = '  -1.41%  '
$ws.Range("D17").Value = "'" + '61.294.23'
$ws.Range("E17").Value = '  -4.24%  '
$ws.Range("D18").Value = "'" + '2.994.43'
$ws.Range("E18").Value = '  -3.46%  '
$ws.Range("D19").Value = "'" + '6.62'
$ws.Range("E19").Value = '  -1.61%  '
$ws.Range("D20").Value = "'" + '464.53'
$ws.Range("E20").Value = '  -5.23%  '
$ws.Range("D21").Value = "'" + '13.20'
$ws.Range("E21").Value = '  -2.45%  '
$ws.Range("D22").Value = "'" + '0.676'
$ws.Range("E22").Value = '  -4.00%  '
$ws.Range("E23").Value = '  -3.73%  '
$ws.Range("D24").Value = "'" + '80.02'
$ws.Range("E24").Value = '  +0.04%  '
$ws.Range("E25").Value = '  -2.51%  '
$ws.Range("D26").Value = "'" + '0.999'
$ws.Range("E26").Value = '  -0.28%  '
$ws.Range("D27").Value = "'" + '2.67'
$ws.Range("E27").Value = '  -2.37%  '
$ws.Range("D28").Value = "'" + '7.79'
$ws.Range("E28").Value = '  -6.90%  '
$ws.Range("D29").Value = "'" + '1.00'
$ws.Range("E29").Value = '  +0.15%  '
$ws.Range("D30").Value = "'" + '1.88'
$ws.Range("E30").Value = '  -1.71%  '
$ws.Range("E31").Value = '  +2.35%  '
$ws.Range("D32").Value = "'" + '25.55'
$ws.Range("E32").Value = '  -2.95%  '
$ws.Range("B33").Value = 'NEARProtocol'
$ws.Range("C33").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D33").Value = "'" + '5.47'
$ws.Range("E33").Value = '  -0.05%  '
$ws.Range("B34").Value = 'OKB'
$ws.Range("C34").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D34").Value = "'" + '55.31'
$ws.Range("E34").Value = '  -3.95%  '
$ws.Range("D35").Value = "'" + '2.27'
$ws.Range("E35").Value = '  -6.27%  '
$ws.Range("E36").Value = '  -3.41%  '
$ws.Range("D37").Value = "'" + '451.31'
$ws.Range("E37").Value = '  -9.06%  '
$ws.Range("D38").Value = "'" + '3.156.36'
$ws.Range("E38").Value = '  -2.93%  '
$ws.Range("D39").Value = "'" + '0.0788'
$ws.Range("E39").Value = '  -1.97%  '
$ws.Range("E40").Value = '  -4.25%  '
$ws.Range("E41").Value = '  +0.07%  '
$ws.Range("D42").Value = "'" + '8.12'
$ws.Range("E42").Value = '  -0.91%  '
$ws.Range("D43").Value = "'" + '2.46'
$ws.Range("E43").Value = '  -8.21%  '
$ws.Range("D44").Value = "'" + '27.18'
$ws.Range("E44").Value = '  +9.35%  '
$ws.Range("E45").Value = '  +0.09%  '
$ws.Range("D46").Value = "'" + '0.244'
$ws.Range("E46").Value = '  -5.62%  '
$ws.Range("D47").Value = "'" + '1.99'
$ws.Range("E47").Value = '  -4.33%  '
$ws.Range("D48").Value = "'" + '119.07'
$ws.Range("E48").Value = '  -2.36%  '
$ws.Range("E49").Value = '  -1.81%  '
$ws.Range("D50").Value = "'" + '0.0₃0495'
$ws.Range("E50").Value = '  -8.34%  '
$ws.Range("D51").Value = "'" + '1.25'
$ws.Range("E51").Value = '  +5.01%  '
